$d = $word.ActiveDocument

# 1) Append "]" to the end of the "Vinicius Kulik Gavioli" paragraph.
$d.Content.Find.Execute(
    "Nome: Vinicius Kulik Gavioli R.A.:819151742", $false, $false, $false,
    $false, $false, $true, 1, $false,
    "Nome: Vinicius Kulik Gavioli R.A.:819151742]", 2)

# 2) Reformat the "Nome : Lucas Bartholetti Palia R.A : 820147772" paragraph
#    (fix the spacing around the colons), keeping its trailing bookmark.
$d.Content.Find.Execute(
    "Nome : Lucas Bartholetti Palia R.A : 820147772", $false, $false, $false,
    $false, $false, $true, 1, $false,
    "Nome: Lucas Bartholetti Palia R.A: 82014772", 2)

# 3) Insert a brand new paragraph ("Marcelo Vinicius...") right before the
#    "Lucas Bartholetti Palia" paragraph, by appending it after the
#    "Vinicius Kulik Gavioli" paragraph that precedes it.
$prevPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Nome: Vinicius Kulik Gavioli*") {
        $prevPara = $p
    }
}
$prevPara.Range.InsertParagraphAfter()
$newPara = $prevPara.Next()
$newPara.Range.Text = "Nome: Marcelo Vinicius Martins da Silva R.A: 820134048"

# 4) Add a new empty paragraph at the very end of the document.
$last = $d.Paragraphs($d.Paragraphs.Count)
$last.Range.InsertParagraphAfter()
